$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 23662
$ws.Cells.Item(2, 4).Value = 34543001
$ws.Cells.Item(3, 3).Value = 59652
$ws.Cells.Item(3, 4).Value = 88299959
$ws.Cells.Item(4, 3).Value = 20241
$ws.Cells.Item(4, 4).Value = 30200738
$ws.Cells.Item(5, 3).Value = 5350
$ws.Cells.Item(5, 4).Value = 8004220
$ws.Cells.Item(6, 3).Value = 1061
$ws.Cells.Item(6, 4).Value = 1589197
$ws.Cells.Item(10, 3).Value = 25344
$ws.Cells.Item(10, 4).Value = 35021449
$ws.Cells.Item(11, 3).Value = 6210
$ws.Cells.Item(11, 4).Value = 9086017
$ws.Cells.Item(12, 3).Value = 17520
$ws.Cells.Item(12, 4).Value = 25915021
$ws.Cells.Item(13, 3).Value = 5453
$ws.Cells.Item(13, 4).Value = 8142636
$ws.Cells.Item(14, 3).Value = 1286
$ws.Cells.Item(14, 4).Value = 1923629
$ws.Cells.Item(15, 3).Value = 233
$ws.Cells.Item(15, 4).Value = 346766
$ws.Cells.Item(17, 3).Value = 6363
$ws.Cells.Item(17, 4).Value = 8623296
$ws.Cells.Item(18, 3).Value = 8666
$ws.Cells.Item(18, 4).Value = 12632523
$ws.Cells.Item(19, 3).Value = 21450
$ws.Cells.Item(19, 4).Value = 31749354
$ws.Cells.Item(20, 3).Value = 6803
$ws.Cells.Item(20, 4).Value = 10167288
$ws.Cells.Item(21, 3).Value = 1634
$ws.Cells.Item(21, 4).Value = 2445802
$ws.Cells.Item(24, 3).Value = 7409
$ws.Cells.Item(24, 4).Value = 10128750
$ws.Cells.Item(25, 3).Value = 4907
$ws.Cells.Item(25, 4).Value = 7158783
$ws.Cells.Item(26, 3).Value = 15119
$ws.Cells.Item(26, 4).Value = 22367007
$ws.Cells.Item(27, 3).Value = 5146
$ws.Cells.Item(27, 4).Value = 7693448
$ws.Cells.Item(28, 3).Value = 1235
$ws.Cells.Item(28, 4).Value = 1851991
$ws.Cells.Item(31, 3).Value = 5232
$ws.Cells.Item(31, 4).Value = 7031990
$ws.Cells.Item(32, 3).Value = 1731
$ws.Cells.Item(32, 4).Value = 2498121
$ws.Cells.Item(33, 3).Value = 4568
$ws.Cells.Item(33, 4).Value = 6711626
$ws.Cells.Item(34, 3).Value = 1851
$ws.Cells.Item(34, 4).Value = 2755922
$ws.Cells.Item(35, 3).Value = 476
$ws.Cells.Item(35, 4).Value = 710541
$ws.Cells.Item(38, 3).Value = 1170
$ws.Cells.Item(38, 4).Value = 1591448
$ws.Cells.Item(39, 3).Value = 11006
$ws.Cells.Item(39, 4).Value = 16050034
$ws.Cells.Item(40, 3).Value = 33947
$ws.Cells.Item(40, 4).Value = 50201835
$ws.Cells.Item(41, 3).Value = 12565
$ws.Cells.Item(41, 4).Value = 18775613
$ws.Cells.Item(42, 3).Value = 3473
$ws.Cells.Item(42, 4).Value = 5199870
$ws.Cells.Item(43, 3).Value = 602
$ws.Cells.Item(43, 4).Value = 901936
$ws.Cells.Item(46, 3).Value = 10355
$ws.Cells.Item(46, 4).Value = 14105688
$ws.Cells.Item(47, 3).Value = 989
$ws.Cells.Item(47, 4).Value = 1431597
$ws.Cells.Item(48, 3).Value = 3722
$ws.Cells.Item(48, 4).Value = 5489209
$ws.Cells.Item(49, 3).Value = 1395
$ws.Cells.Item(49, 4).Value = 2085464
$ws.Cells.Item(50, 3).Value = 429
$ws.Cells.Item(50, 4).Value = 641000
$ws.Cells.Item(52, 3).Value = 2417
$ws.Cells.Item(52, 4).Value = 3364281
$ws.Cells.Item(53, 3).Value = 352
$ws.Cells.Item(53, 4).Value = 511784
$ws.Cells.Item(54, 3).Value = 947
$ws.Cells.Item(54, 4).Value = 1402492
$ws.Cells.Item(55, 3).Value = 384
$ws.Cells.Item(55, 4).Value = 573976
$ws.Cells.Item(56, 3).Value = 131
$ws.Cells.Item(56, 4).Value = 196378
$ws.Cells.Item(58, 3).Value = 458
$ws.Cells.Item(58, 4).Value = 650481
$ws.Cells.Item(59, 3).Value = 9999
$ws.Cells.Item(59, 4).Value = 14529958
$ws.Cells.Item(60, 3).Value = 30438
$ws.Cells.Item(60, 4).Value = 44914535
$ws.Cells.Item(61, 3).Value = 10557
$ws.Cells.Item(61, 4).Value = 15779522
$ws.Cells.Item(62, 3).Value = 2928
$ws.Cells.Item(62, 4).Value = 4382138
$ws.Cells.Item(63, 3).Value = 516
$ws.Cells.Item(63, 4).Value = 773639
$ws.Cells.Item(66, 3).Value = 9893
$ws.Cells.Item(66, 4).Value = 13259323
$ws.Cells.Item(67, 3).Value = 2709
$ws.Cells.Item(67, 4).Value = 3954918
$ws.Cells.Item(68, 3).Value = 7371
$ws.Cells.Item(68, 4).Value = 10878685
$ws.Cells.Item(69, 3).Value = 2612
$ws.Cells.Item(69, 4).Value = 3903113
$ws.Cells.Item(70, 3).Value = 856
$ws.Cells.Item(70, 4).Value = 1282049
$ws.Cells.Item(71, 3).Value = 174
$ws.Cells.Item(71, 4).Value = 259612
$ws.Cells.Item(73, 3).Value = 2847
$ws.Cells.Item(73, 4).Value = 3878333
$ws.Cells.Item(74, 3).Value = 874
$ws.Cells.Item(74, 4).Value = 1267893
$ws.Cells.Item(75, 3).Value = 2994
$ws.Cells.Item(75, 4).Value = 4425546
$ws.Cells.Item(76, 3).Value = 1189
$ws.Cells.Item(76, 4).Value = 1779659
$ws.Cells.Item(77, 3).Value = 412
$ws.Cells.Item(77, 4).Value = 618000
$ws.Cells.Item(80, 3).Value = 1767
$ws.Cells.Item(80, 4).Value = 2371711
$ws.Cells.Item(86, 3).Value = 7022
$ws.Cells.Item(86, 4).Value = 10271210
$ws.Cells.Item(87, 3).Value = 20193
$ws.Cells.Item(87, 4).Value = 29875547
$ws.Cells.Item(88, 3).Value = 6627
$ws.Cells.Item(88, 4).Value = 9905215
$ws.Cells.Item(89, 3).Value = 1752
$ws.Cells.Item(89, 4).Value = 2623655
$ws.Cells.Item(93, 3).Value = 6299
$ws.Cells.Item(93, 4).Value = 8494577
$ws.Cells.Item(94, 3).Value = 19297
$ws.Cells.Item(94, 4).Value = 28028945
$ws.Cells.Item(95, 3).Value = 44795
$ws.Cells.Item(95, 4).Value = 66086048
$ws.Cells.Item(96, 3).Value = 14321
$ws.Cells.Item(96, 4).Value = 21382995
$ws.Cells.Item(97, 3).Value = 3828
$ws.Cells.Item(97, 4).Value = 5727884
$ws.Cells.Item(98, 3).Value = 653
$ws.Cells.Item(98, 4).Value = 977862
$ws.Cells.Item(99, 3).Value = 32
$ws.Cells.Item(99, 4).Value = 46164
$ws.Cells.Item(101, 3).Value = 16501
$ws.Cells.Item(101, 4).Value = 22425353
$ws.Cells.Item(102, 3).Value = 22076
$ws.Cells.Item(102, 4).Value = 32102834
$ws.Cells.Item(103, 3).Value = 49925
$ws.Cells.Item(103, 4).Value = 73568262
$ws.Cells.Item(104, 3).Value = 15591
$ws.Cells.Item(104, 4).Value = 23252519
$ws.Cells.Item(105, 3).Value = 3994
$ws.Cells.Item(105, 4).Value = 5967314
$ws.Cells.Item(106, 3).Value = 646
$ws.Cells.Item(106, 4).Value = 966054
$ws.Cells.Item(109, 3).Value = 19553
$ws.Cells.Item(109, 4).Value = 26386817
$ws.Cells.Item(110, 3).Value = 8587
$ws.Cells.Item(110, 4).Value = 12545461
$ws.Cells.Item(111, 3).Value = 22246
$ws.Cells.Item(111, 4).Value = 32929074
$ws.Cells.Item(112, 3).Value = 7715
$ws.Cells.Item(112, 4).Value = 11518314
$ws.Cells.Item(113, 3).Value = 1869
$ws.Cells.Item(113, 4).Value = 2796094
$ws.Cells.Item(117, 3).Value = 7003
$ws.Cells.Item(117, 4).Value = 9535160
$ws.Cells.Item(118, 3).Value = 21483
$ws.Cells.Item(118, 4).Value = 31239977
$ws.Cells.Item(119, 3).Value = 52955
$ws.Cells.Item(119, 4).Value = 78112377
$ws.Cells.Item(120, 3).Value = 15943
$ws.Cells.Item(120, 4).Value = 23810363
$ws.Cells.Item(121, 3).Value = 3959
$ws.Cells.Item(121, 4).Value = 5921997
$ws.Cells.Item(122, 3).Value = 808
$ws.Cells.Item(122, 4).Value = 1210212
$ws.Cells.Item(124, 3).Value = 7
$ws.Cells.Item(124, 4).Value = 10500
$ws.Cells.Item(125, 3).Value = 18270
$ws.Cells.Item(125, 4).Value = 25113464
$ws.Cells.Item(126, 3).Value = 29462
$ws.Cells.Item(126, 4).Value = 43155788
$ws.Cells.Item(127, 3).Value = 88631
$ws.Cells.Item(127, 4).Value = 131293383
$ws.Cells.Item(128, 3).Value = 39350
$ws.Cells.Item(128, 4).Value = 58816398
$ws.Cells.Item(129, 3).Value = 12244
$ws.Cells.Item(129, 4).Value = 18339357
$ws.Cells.Item(130, 3).Value = 2491
$ws.Cells.Item(130, 4).Value = 3730909
$ws.Cells.Item(134, 3).Value = 28980
$ws.Cells.Item(134, 4).Value = 40419381
